# "remove MR from main, utils"
#
# 1. "Propellant Combinations" sheet: no longer the active tab; selection
#    changes to D1:G4 (active cell D1).
# 2. "Continuous Inputs" sheet: insert a new column C ("Exit pressure (psi)")
#    between "Chamber pressure (psi)" (B) and "Thrust-to-Weight ratio" (old C,
#    now D), with values 7 / 13 / 0.5 for the three rows. This sheet becomes
#    the active tab, selection C7.

$wb = $excel.ActiveWorkbook

# ---- Propellant Combinations: update selection, no longer active tab ----
$wsProp = $wb.Worksheets.Item("Propellant Combinations")
$wsProp.Activate()
$wsProp.Range("D1:G4").Select()

# ---- Continuous Inputs: insert "Exit pressure (psi)" column ----
$wsCont = $wb.Worksheets.Item("Continuous Inputs")
$wsCont.Activate()

$wsCont.Columns("C").Insert()
$wsCont.Columns("C").ColumnWidth = $wsCont.Columns("B").ColumnWidth

$wsCont.Range("C1").Value = "Exit pressure (psi)"
$wsCont.Range("C2").Value = 7
$wsCont.Range("C3").Value = 13
$wsCont.Range("C4").Value = 0.5

# Match the number format of the (now-shifted) "Thrust-to-Weight ratio"
# column so the new C4 cell picks up the same style as D4.
$wsCont.Range("C4").NumberFormat = $wsCont.Range("D4").NumberFormat

$wsCont.Range("C7").Select()
